$wb = $excel.ActiveWorkbook

# Data to insert: sheet name -> [date, price] inserted as a new row just
# before the existing last row (2025-06-25), pushing it down.
$updates = @(
    @{ Sheet = "Gaz"; Date = "2025-06-24"; Price = 40.9 },
    @{ Sheet = "CO2"; Date = "2025-06-24"; Price = 71.88 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)

    # Insert a new blank row at row 10, shifting the existing row 10
    # (2025-06-25) down to row 11.
    $ws.Rows.Item(10).Insert()

    # Force column A to text so the date string isn't auto-converted into a
    # date serial number (the source sheet stores dates as plain text),
    # then restore the default "Normal" style so no stray number-format
    # style lingers on the cell.
    $ws.Cells.Item(10, 1).NumberFormat = "@"
    $ws.Cells.Item(10, 1).Value = $u.Date
    $ws.Cells.Item(10, 1).Style = "Normal"

    # Fill in the newly inserted row's price.
    $ws.Cells.Item(10, 2).Value = $u.Price
}
